$d = $word.ActiveDocument

# --- 1. Insert the six new paragraphs right after "Activités / Tâches :" ---
# Use Find & Replace with paragraph-mark codes (^p) so every new line lands
# as its own <w:p>, inheriting the Titre3 style of the matched paragraph to
# start with; we fix up the individual paragraph styles/numbering below.
$newLines = @(
    "Création d’un compte sur le dépôt GitHub.",
    "Durée : 15 min",
    "Liaison avec le projet de groupe.",
    "Durée : 5 min",
    "Clonage du dépôt en local et enrichissement en données de celui-ci (Gantt).",
    "Durée : 20 min"
)
$replacement = "Activités / Tâches :^p" + ($newLines -join "^p")

$findRange = $d.Content
$found = $findRange.Find.Execute("Activités / Tâches :", $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

# --- 2. Locate the anchor paragraph ("Activités / Tâches :") again, now
#        that the new paragraphs exist right after it. ---
$paragraphs = $d.Paragraphs
$anchorIndex = -1
$ganttTemplate = $null
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    $text = $para.Range.Text
    if ($anchorIndex -eq -1 -and $text -like "Activit*T?ches*") {
        $anchorIndex = $i
    }
    if ($ganttTemplate -eq $null -and $text -like "Finalisation de la planification*") {
        # Grab the list template already used by the existing Gantt bullet
        # (numId 3) so the new bullet paragraphs continue the same list
        # instead of minting a brand-new numbering definition.
        $ganttTemplate = $para.Range.ListFormat.ListTemplate
    }
}

# --- 3. Re-style the six freshly inserted paragraphs. ---
$styles = @(
    "Paragraphedeliste",
    "Titre4",
    "Paragraphedeliste",
    "Titre4",
    "Paragraphedeliste",
    "Titre4"
)

for ($j = 0; $j -lt $styles.Length; $j++) {
    $p = $d.Paragraphs.Item($anchorIndex + 1 + $j)
    $p.Style = $styles[$j]
    if ($styles[$j] -eq "Paragraphedeliste") {
        $p.Range.ListFormat.ApplyListTemplateWithLevel($ganttTemplate, $true, 1, $false, 1)
    }
}
